$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old B2:C3 data that is no longer part of the table
# (the "Data" / "Pedidos no Dia" columns are no longer filled per-row)
$ws.Range("B2:C3").ClearContents()

# Column A: "Tarefa" names, Column D: "Observação" classification
$rows = @(
    @("NEONATURE",     "BAIXA DEMANDA / SUPORTE MEDIO"),
    @("GWS",           "ALTA DEMANDA / SUPORTE MEDIO"),
    @("VANGUARDA",     "MEDIA DEMANDA / SUPORTE ALTO"),
    @("ASTROMIC",      "ALTA DEMANDA / SUPORTE ALTO"),
    @("NEURO BETES",   "POUCA DEMANDA / SUPORTE MEDIO"),
    @("HERA",          "MEDIA DEMANDA / SUPORTE MEDIO"),
    @("LUNO",          "ALTA DEMANDA / SUPORTE ALTO"),
    @("MF",            "BAIXÍSSIMA DEMANDA / SUPORTE BAIXÍSSIMO EXCETO NOS DIAS QUE TEM PEDIDO"),
    @("RADT",          "BAIXA DEMANDA / SUPORTE BAIXO"),
    @("DESAGITA",      "MEDIA DEMANDA / SUPORTE ALTO"),
    @("DIABETINA",     "0 DEMANDA / SUPORTE NULO"),
    @("SAUDE FITNESS", "0 DEMANDA / SUPORTE NULO")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $r = $r + 1
}

Write-Host "done"
